$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.114.32'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '1.599.65'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.07'
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.483'
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0615'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.30'
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +2.49%  '
$ws.Range("D12").Value = '1.820.61'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = '1.600.24'
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").Value = '26.095.97'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.93'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '204.04'
$ws.Range("E20").Value = '  +4.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.30'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.04'
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.92'
$ws.Range("E24").Value = '  +13.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.64'
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  -7.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.24'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0479'
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.15'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.91'
$ws.Range("E33").Value = '  -3.70%  '
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").Value = '1.133.11'
$ws.Range("E36").Value = '  +3.39%  '
$ws.Range("E37").Value = '  +7.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.800'
$ws.Range("E38").Value = '  +2.49%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").Value = '1.733.89'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.20'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.26'
$ws.Range("E47").Value = '  +2.01%  '
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₇0972'
$ws.Range("E49").Value = '  -12.13%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.406'
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  +0.14%  '
